# "corrected 2023 start date"
# The hunting season start date for 2023 (row 6, column B) was recorded as
# 2023-10-23 (serial 45222) but should have been 2023-10-21 (serial 45220).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 45220

# Leave the selection where it was when the workbook was last saved.
$ws.Range("E9").Select()
